$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2997
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 4994
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 4994
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -5344
$ws.Range("H55").Value = 2776.8
$ws.Range("I55").Value = 306.66666
$ws.Range("J55").Value = 3835.4285
$ws.Range("K55").Value = 306.66666
$ws.Range("L55").Value = 3835.4285
$ws.Range("M55").Value = -92.66665999999998
$ws.Range("N55").Value = -4263.4285
$ws.Range("H96").Value = 2348.3845
$ws.Range("J96").Value = 2694.5
$ws.Range("L96").Value = 8083.5
$ws.Range("N96").Value = -10829.5
$ws.Range("H100").Value = 3332.5334
$ws.Range("I100").Value = 1996.3334
$ws.Range("J100").Value = 3666.5833
$ws.Range("K100").Value = 1996.3334
$ws.Range("L100").Value = 3666.5833
$ws.Range("M100").Value = -1455.3334
$ws.Range("N100").Value = -4748.5833
$ws.Range("H113").Value = 71432070
$ws.Range("I113").Value = 33335664
$ws.Range("J113").Value = 100004376
$ws.Range("K113").Value = 33335664
$ws.Range("L113").Value = 100004376
$ws.Range("M113").Value = -33332410
$ws.Range("N113").Value = -100010884
$ws.Range("H114").Value = 54887.5
$ws.Range("J114").Value = 54887.5
$ws.Range("L114").Value = 54887.5
$ws.Range("N114").Value = -63565.5
$ws.Range("H138").Value = 2903.9788
$ws.Range("I138").Value = 1360.8125
$ws.Range("J138").Value = 3700.4517
$ws.Range("K138").Value = 4082.4375
$ws.Range("L138").Value = 11101.3551
$ws.Range("M138").Value = 1057.5625
$ws.Range("N138").Value = -21381.3551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 458.7
$ws.Range("I4").Value = 453.8889
$ws.Range("J4").Value = 502
$ws.Range("K4").Value = 453.8889
$ws.Range("L4").Value = 502
$ws.Range("M4").Value = -337.8889
$ws.Range("N4").Value = -734
$ws.Range("H11").Value = 5062502
$ws.Range("I11").Value = 20000000
$ws.Range("K11").Value = 20000000
$ws.Range("M11").Value = -19999856
$ws.Range("H110").Value = 2447.2727
$ws.Range("I110").Value = 2250.889
$ws.Range("J110").Value = 3331
$ws.Range("K110").Value = 2250.889
$ws.Range("L110").Value = 3331
$ws.Range("M110").Value = -205.8890000000001
$ws.Range("N110").Value = -7421
$ws.Range("H122").Value = 2716.077
$ws.Range("I122").Value = 2550
$ws.Range("J122").Value = 2858.4285
$ws.Range("K122").Value = 7650
$ws.Range("L122").Value = 8575.2855
$ws.Range("M122").Value = -5200
$ws.Range("N122").Value = -13475.2855
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3449.75
$ws.Range("I86").Value = 3933
$ws.Range("K86").Value = 3933
$ws.Range("M86").Value = -2810
$ws.Range("H89").Value = 3449.75
$ws.Range("I89").Value = 3933
$ws.Range("K89").Value = 19665
$ws.Range("M89").Value = -14049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 627.5
$ws.Range("I3").Value = 170
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 170
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -57
$ws.Range("N3").Value = -2226
$ws.Range("H31").Value = 868428.4
$ws.Range("I31").Value = 2840.7273
$ws.Range("J31").Value = 1463519.9
$ws.Range("K31").Value = 2840.7273
$ws.Range("L31").Value = 1463519.9
$ws.Range("M31").Value = -2545.7273
$ws.Range("N31").Value = -1464109.9
$ws.Range("H34").Value = 868428.4
$ws.Range("I34").Value = 2840.7273
$ws.Range("J34").Value = 1463519.9
$ws.Range("K34").Value = 2840.7273
$ws.Range("L34").Value = 1463519.9
$ws.Range("M34").Value = -2638.7273
$ws.Range("N34").Value = -1463923.9
$ws.Range("H99").Value = 4071.3333
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 4107
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 4107
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -7103
$ws.Range("H107").Value = 2027.3334
$ws.Range("I107").Value = 908.875
$ws.Range("K107").Value = 908.875
$ws.Range("M107").Value = 1011.125
$ws.Range("H126").Value = 4071.3333
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 4107
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 12321
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -17261

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 98852.42999999999
$ws.Range("J9").Value = 113661.336
$ws.Range("L9").Value = 340984.008
$ws.Range("N9").Value = -341432.008
$ws.Range("H11").Value = 33.5
$ws.Range("I11").Value = 33.5
$ws.Range("K11").Value = 100.5
$ws.Range("M11").Value = 39.5
$ws.Range("H37").Value = 76899.8
$ws.Range("J37").Value = 76899.8
$ws.Range("L37").Value = 230699.4
$ws.Range("N37").Value = -230923.4
$ws.Range("H54").Value = 9000.571
$ws.Range("J54").Value = 11000
$ws.Range("L54").Value = 33000
$ws.Range("N54").Value = -34118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6065.5835
$ws.Range("I70").Value = 4814.8335
$ws.Range("K70").Value = 4814.8335
$ws.Range("M70").Value = -4544.8335
$ws.Range("H73").Value = 6065.5835
$ws.Range("I73").Value = 4814.8335
$ws.Range("K73").Value = 4814.8335
$ws.Range("M73").Value = -3878.8335
$ws.Range("H132").Value = 30310854
$ws.Range("I132").Value = 40004876
$ws.Range("J132").Value = 17031.875
$ws.Range("K132").Value = 120014628
$ws.Range("L132").Value = 51095.625
$ws.Range("M132").Value = -120012098
$ws.Range("N132").Value = -56155.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 4765.5
$ws.Range("I35").Value = 4765.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4765.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4429.5
$ws.Range("N35").ClearContents()
$ws.Range("H41").Value = 32495
$ws.Range("J41").Value = 32495
$ws.Range("L41").Value = 32495
$ws.Range("N41").Value = -33371
$ws.Range("H132").Value = 786296.1
$ws.Range("I132").Value = 24855.625
$ws.Range("J132").Value = 2004601
$ws.Range("K132").Value = 74566.875
$ws.Range("L132").Value = 6013803
$ws.Range("M132").Value = -72036.875
$ws.Range("N132").Value = -6018863
